{"js": "// Replace the three-digit-by-one-digit multiplication expressions in the\n// document's table cells with the new values from the commit.\nconst replacements = [\n  [\"402\u00d74=\", \"126\u00d75=\"],\n  [\"347\u00d72=\", \"776\u00d72=\"],\n  [\"800\u00d76=\", \"527\u00d73=\"],\n  [\"295\u00d74=\", \"953\u00d79=\"],\n  [\"661\u00d74=\", \"528\u00d73=\"],\n  [\"308\u00d72=\", \"748\u00d73=\"],\n  [\"921\u00d75=\", \"756\u00d77=\"],\n  [\"923\u00d79=\", \"610\u00d77=\"],\n  [\"543\u00d75=\", \"532\u00d74=\"],\n  [\"680\u00d75=\", \"499\u00d78=\"],\n  [\"914\u00d77=\", \"303\u00d75=\"],\n  [\"242\u00d78=\", \"289\u00d76=\"],\n  [\"512\u00d72=\", \"370\u00d74=\"],\n  [\"902\u00d72=\", \"525\u00d72=\"],\n  [\"618\u00d74=\", \"833\u00d78=\"],\n  [\"355\u00d72=\", \"272\u00d74=\"],\n  [\"787\u00d74=\", \"406\u00d72=\"],\n  [\"438\u00d79=\", \"659\u00d73=\"],\n  [\"291\u00d73=\", \"507\u00d79=\"],\n  [\"624\u00d76=\", \"298\u00d75=\"],\n  [\"714\u00d75=\", \"323\u00d77=\"],\n  [\"163\u00d78=\", \"624\u00d78=\"],\n  [\"613\u00d73=\", \"198\u00d76=\"],\n  [\"692\u00d72=\", \"655\u00d77=\"],\n  [\"445\u00d78=\", \"518\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication expressions in the\n# document's table cells with the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"402\u00d74=\", \"126\u00d75=\"),\n    @(\"347\u00d72=\", \"776\u00d72=\"),\n    @(\"800\u00d76=\", \"527\u00d73=\"),\n    @(\"295\u00d74=\", \"953\u00d79=\"),\n    @(\"661\u00d74=\", \"528\u00d73=\"),\n    @(\"308\u00d72=\", \"748\u00d73=\"),\n    @(\"921\u00d75=\", \"756\u00d77=\"),\n    @(\"923\u00d79=\", \"610\u00d77=\"),\n    @(\"543\u00d75=\", \"532\u00d74=\"),\n    @(\"680\u00d75=\", \"499\u00d78=\"),\n    @(\"914\u00d77=\", \"303\u00d75=\"),\n    @(\"242\u00d78=\", \"289\u00d76=\"),\n    @(\"512\u00d72=\", \"370\u00d74=\"),\n    @(\"902\u00d72=\", \"525\u00d72=\"),\n    @(\"618\u00d74=\", \"833\u00d78=\"),\n    @(\"355\u00d72=\", \"272\u00d74=\"),\n    @(\"787\u00d74=\", \"406\u00d72=\"),\n    @(\"438\u00d79=\", \"659\u00d73=\"),\n    @(\"291\u00d73=\", \"507\u00d79=\"),\n    @(\"624\u00d76=\", \"298\u00d75=\"),\n    @(\"714\u00d75=\", \"323\u00d77=\"),\n    @(\"163\u00d78=\", \"624\u00d78=\"),\n    @(\"613\u00d73=\", \"198\u00d76=\"),\n    @(\"692\u00d72=\", \"655\u00d77=\"),\n    @(\"445\u00d78=\", \"518\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n}\n"}
